$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "helix jump 2" row (row 5). Rows below shift up by one,
# so the former row 6 ("stretchy taxi") becomes the new row 5, etc.
$ws.Rows.Item(5).Delete() | Out-Null

# Append a new last row (12) for "taxi game" / "com.singleton.strechy",
# re-using the formatting of the row above it so the new cells keep the
# same style as the rest of the table.
$ws.Range("A11:B11").Copy() | Out-Null
$ws.Range("A12:B12").PasteSpecial(-4122) | Out-Null
$ws.Range("A12").Value = "taxi game"
$ws.Range("B12").Value = "com.singleton.strechy"

# Update the active selection to A5, matching the saved view state.
$ws.Range("A5").Select() | Out-Null
